$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<do>"
$ws.Range("C2").Value = 30

$ws.Range("B3").Value = "<ki>"
$ws.Range("C3").Value = 31

$ws.Range("B4").Value = "<now>"
$ws.Range("C4").Value = 34

$ws.Range("C5").Value = 25

$ws.Range("B6").Value = "<not>"
$ws.Range("C6").Value = 24

$ws.Range("B7").Value = "<otha>"
$ws.Range("C7").Value = 28

$ws.Range("B8").Value = "<four>"

$ws.Range("C9").Value = 35

$ws.Range("C10").Value = 22

$ws.Range("B11").Value = "<ben>"
$ws.Range("C11").Value = 28

$ws.Range("C12").Value = 20

$ws.Range("C13").Value = 31

$ws.Range("B14").Value = "<they>"
$ws.Range("C14").Value = 33

$ws.Range("B15").Value = "<in>"
$ws.Range("C15").Value = 28

$ws.Range("B16").Value = "<pace>"
$ws.Range("C16").Value = 34

$ws.Range("C17").Value = 29

$ws.Range("C18").Value = 27
